# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund-holding detail) right before the
#   "总计" (totals) sheet, using the same header/style layout as "2021-Q4".
# - Add a new leading row to the "总计" sheet summarising the 2022-Q1 data,
#   pushing the existing 2021-Q4 / 2021-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell without Excel's "looks like a
# number" auto-conversion (e.g. "020009" -> 20009, "7.10" -> 7.1) and without
# leaving a stray quote-prefix / number-format style behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $blank = $range.Worksheet.Range("ZZ1000")
    $range.NumberFormat = "@"
    $range.Value = $text
    $blank.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats - restores the default (no-op) style
}

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet by duplicating "2021-Q4" (identical header
#    text/style layout), placing it right after "2021-Q4" i.e. right before
#    "总计".
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# extend the existing data-row formatting (A2:H2) down to rows 3-4
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A3:H4").PasteSpecial(-4122)

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "011466"
$newSheet.Range("C2").Value = "兴业医疗保健混合A"
Set-TextValue $newSheet.Range("D2") "7.10"
Set-TextValue $newSheet.Range("E2") "84.68"
Set-TextValue $newSheet.Range("F2") "3.32"
Set-TextValue $newSheet.Range("G2") "0.2357"
$newSheet.Range("H2").Value = 10

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "020009"
$newSheet.Range("C3").Value = "国泰金鹏蓝筹混合"
Set-TextValue $newSheet.Range("D3") "5.45"
Set-TextValue $newSheet.Range("E3") "89.63"
Set-TextValue $newSheet.Range("F3") "2.37"
Set-TextValue $newSheet.Range("G3") "0.1292"
$newSheet.Range("H3").Value = 8

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "011467"
$newSheet.Range("C4").Value = "兴业医疗保健混合C"
Set-TextValue $newSheet.Range("D4") "2.22"
Set-TextValue $newSheet.Range("E4") "84.68"
Set-TextValue $newSheet.Range("F4") "3.32"
Set-TextValue $newSheet.Range("G4") "0.0737"
$newSheet.Range("H4").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: push the two existing data rows down by one and
#    insert the new 2022-Q1 summary row at the top of the data.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# extend the existing data-row formatting (A2:D2) down to rows 3-4
$total.Range("A2:D2").Copy()
$total.Range("A3:D4").PasteSpecial(-4122)

# Row 2 (new)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.44

# Row 3 (was row 2: 2021-Q4)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.17

# Row 4 (was row 3: 2021-Q2)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.18
